$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(7, 2).Range.Text = "CAD MOV SON "
$t.Cell(8, 2).Range.Text = "CUT "
$t.Cell(9, 2).Range.Text = "CAD SON "
$t.Cell(10, 2).Range.Text = "CAD SON "
$t.Cell(11, 2).Range.Text = "CAD MOV SON "
$t.Cell(12, 2).Range.Text = "CAD REG SON "
$t.Cell(13, 2).Range.Text = "CUT SON "
$t.Cell(14, 2).Range.Text = "CUT SON "
$t.Cell(15, 2).Range.Text = "CUT SON "
$t.Cell(16, 2).Range.Text = "CAD SON "
$t.Cell(17, 2).Range.Text = "CAD SON "
$t.Cell(18, 2).Range.Text = "CAD MOV SON "
$t.Cell(19, 2).Range.Text = "C3 CAD MOV SON "
$t.Cell(20, 2).Range.Text = "C3 CAD SON "
$t.Cell(21, 2).Range.Text = "C3 CAD SON "
$t.Cell(22, 2).Range.Text = "CAD SON "
$t.Cell(23, 2).Range.Text = "CAD SON "
$t.Cell(24, 2).Range.Text = "CAD SON "
$t.Cell(25, 2).Range.Text = "CUT SON "
$t.Cell(26, 2).Range.Text = "CAD MOV SON "
$t.Cell(27, 2).Range.Text = "CAD SON "
$t.Cell(28, 2).Range.Text = "CAD SON "
$t.Cell(29, 2).Range.Text = "CAD SON "
$t.Cell(30, 2).Range.Text = "CAD SON "
$t.Cell(31, 2).Range.Text = "CAD SON "
$t.Cell(32, 2).Range.Text = "CAD SON "
$t.Cell(33, 2).Range.Text = "C3 CAD SON "
$t.Cell(34, 2).Range.Text = "CAD SON "
$t.Cell(35, 2).Range.Text = "CAD SON "
$t.Cell(36, 2).Range.Text = "C3 CAD SON "
$t.Cell(37, 2).Range.Text = "CAD SON "
$t.Cell(38, 2).Range.Text = "CUT SON "
$t.Cell(39, 2).Range.Text = "CAD SON "
$t.Cell(40, 2).Range.Text = "CUT SON"
$t.Cell(41, 2).Range.Text = "CAD SON "
$t.Cell(42, 2).Range.Text = "CAD MOV SON "
$t.Cell(43, 2).Range.Text = "CAD SON "
$t.Cell(44, 2).Range.Text = "CAD SON "
$t.Cell(45, 2).Range.Text = "CAD SON "
$t.Cell(46, 2).Range.Text = "CAD SON "
$t.Cell(47, 2).Range.Text = "CAD SON "
$t.Cell(48, 2).Range.Text = "CAD MOV SON "
$t.Cell(49, 2).Range.Text = "CAD MET SON "
$t.Cell(50, 2).Range.Text = "CAD SON "
$t.Cell(51, 2).Range.Text = "CAD SON "
$t.Cell(52, 2).Range.Text = "CAD SON"
$t.Cell(53, 2).Range.Text = "CUT MET SON "
$t.Cell(54, 2).Range.Text = "CAD SON "
$t.Cell(55, 2).Range.Text = "CAD SON "
$t.Cell(56, 2).Range.Text = "CAD SON "
$t.Cell(57, 2).Range.Text = "CAD SON "
$t.Cell(58, 2).Range.Text = "C3 CAD SON "
$t.Cell(59, 2).Range.Text = "CAD SON "
$t.Cell(60, 2).Range.Text = "CAD SON "
$t.Cell(61, 2).Range.Text = "CAD SON "
$t.Cell(62, 2).Range.Text = "CAD SON "
$t.Cell(63, 2).Range.Text = "CAD MET SON "
$t.Cell(64, 2).Range.Text = "CAD SON "
$t.Cell(65, 2).Range.Text = "CUT "
$t.Cell(66, 2).Range.Text = "CAD SON "
$t.Cell(67, 2).Range.Text = "CAD SON"
$t.Cell(68, 2).Range.Text = "CAD SON "
$t.Cell(69, 2).Range.Text = "CAD SON"
$t.Cell(70, 2).Range.Text = "CUT SON "
$t.Cell(71, 2).Range.Text = "CAD SON"
$t.Cell(72, 2).Range.Text = "C3 CAD SON"
$t.Cell(73, 2).Range.Text = "CAD SON"
$t.Cell(74, 2).Range.Text = "CAD MET MOV SON "
$t.Cell(75, 2).Range.Text = "CAD MOV SON "
$t.Cell(76, 2).Range.Text = "CAD SON"
$t.Cell(77, 2).Range.Text = "CAD SON"
$t.Cell(78, 2).Range.Text = "CAD SON"
$t.Cell(79, 2).Range.Text = "CAD SON"
$t.Cell(80, 2).Range.Text = "CAD SON"
$t.Cell(81, 2).Range.Text = "CAD "
$t.Cell(82, 2).Range.Text = "CAD SON"
$t.Cell(83, 2).Range.Text = "CAD SON"
$t.Cell(84, 2).Range.Text = "CAD SON"
$t.Cell(85, 2).Range.Text = "CAD SON"
$t.Cell(86, 2).Range.Text = "CAD SON"
$t.Cell(87, 2).Range.Text = "CAD SON"
$t.Cell(88, 2).Range.Text = "CUT SON "
$t.Cell(89, 2).Range.Text = "CAD SON "
$t.Cell(90, 2).Range.Text = "CAD REG SON "
$t.Cell(91, 2).Range.Text = "CAD SON"
$t.Cell(92, 2).Range.Text = "CAD SON"
$t.Cell(93, 2).Range.Text = "CAD SON "
$t.Cell(94, 2).Range.Text = "CAD SON"
$t.Cell(95, 2).Range.Text = "SON"
$t.Cell(96, 2).Range.Text = "SON"
$t.Cell(97, 2).Range.Text = "SON "
$t.Cell(98, 2).Range.Text = "SON "

# Add an extra empty paragraph after the table (before the final paragraph / sectPr)
$insertPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$insertPoint.InsertParagraphBefore()

